$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Range("E$row").Value = 3.67
    $ws.Range("G$row").Value = 3
}
